# "Corr results with infer no filter"
# Updates correlation-analysis numbers on the "all_tools" and "infer" sheets.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet: all_tools
# ---------------------------------------------------------------------
$wsAll = $wb.Worksheets.Item("all_tools")

$wsAll.Range("G10").Value = 822
$wsAll.Range("G11").Value = 822
$wsAll.Range("G12").Value = 822

$wsAll.Range("G25").Value = 40
$wsAll.Range("I25").Value = -0.1807753815155468
$wsAll.Range("J25").Value = 0.3541954904764164
$wsAll.Range("K25").Value = -0.2576049186596542
$wsAll.Range("L25").Value = 0.3354345184285685

$wsAll.Range("G26").Value = 40
$wsAll.Range("I26").Value = -0.1807753815155468
$wsAll.Range("J26").Value = 0.3541954904764164
$wsAll.Range("K26").Value = -0.2666436877354316
$wsAll.Range("L26").Value = 0.3181414648703181

$wsAll.Range("G27").Value = 40
$wsAll.Range("I27").Value = 0.3539900381483285
$wsAll.Range("J27").Value = 0.07056136851892029
$wsAll.Range("K27").Value = 0.4341802833034056
$wsAll.Range("L27").Value = 0.09288178063084394

$wsAll.Range("G28").Value = 40
$wsAll.Range("K28").Value = -0.2493004677260264
$wsAll.Range("L28").Value = 0.3517858440384553

$wsAll.Range("G29").Value = 40
$wsAll.Range("K29").Value = -0.1491396897503261
$wsAll.Range("L29").Value = 0.5814513259975999

# ---------------------------------------------------------------------
# Sheet: infer
# ---------------------------------------------------------------------
$wsInfer = $wb.Worksheets.Item("infer")

# Column K (11) widened slightly to match the other sheets.
$wsInfer.Columns.Item(11).ColumnWidth = 20.9

$wsInfer.Range("F10").Value = 23
$wsInfer.Range("G10").Value = 24
$wsInfer.Range("I10").Value = -0.004001088444105332
$wsInfer.Range("J10").Value = 0.9739374982488735
$wsInfer.Range("K10").Value = -0.004014544573191041
$wsInfer.Range("L10").Value = 0.9779259755977681

$wsInfer.Range("F11").Value = 23
$wsInfer.Range("G11").Value = 24
$wsInfer.Range("I11").Value = -0.01888959038201153
$wsInfer.Range("J11").Value = 0.8717029852220165
$wsInfer.Range("K11").Value = -0.02206154843492711
$wsInfer.Range("L11").Value = 0.8791305588586364

$wsInfer.Range("F12").Value = 23
$wsInfer.Range("G12").Value = 24
$wsInfer.Range("I12").Value = -0.2476064162497625
$wsInfer.Range("J12").Value = 0.0331142128596994
$wsInfer.Range("K12").Value = -0.3099127537051835
$wsInfer.Range("L12").Value = 0.02850794708171081

$wsInfer.Range("F25").Value = 1
$wsInfer.Range("G25").Value = 1
$wsInfer.Range("I25").Value = -0.3535533905932737
$wsInfer.Range("J25").Value = 0.1037416782365415
$wsInfer.Range("K25").Value = -0.4200840252084029
$wsInfer.Range("L25").Value = 0.105228057983522

$wsInfer.Range("F26").Value = 1
$wsInfer.Range("G26").Value = 1
$wsInfer.Range("I26").Value = -0.1649915822768611
$wsInfer.Range("J26").Value = 0.4476990724652935
$wsInfer.Range("K26").Value = -0.1960392117639214
$wsInfer.Range("L26").Value = 0.4668248490265503

$wsInfer.Range("F27").Value = 1
$wsInfer.Range("G27").Value = 1
$wsInfer.Range("I27").Value = 0.02366905341655754
$wsInfer.Range("J27").Value = 0.9135633303377861
$wsInfer.Range("K27").Value = 0.02802621677476181
$wsInfer.Range("L27").Value = 0.9179387985999929

$wsInfer.Range("F28").Value = 1
$wsInfer.Range("G28").Value = 1
$wsInfer.Range("I28").Value = -0.2625754538144587
$wsInfer.Range("J28").Value = 0.2314460271038938
$wsInfer.Range("K28").Value = -0.3089716991054783
$wsInfer.Range("L28").Value = 0.2442606266224961

$wsInfer.Range("F29").Value = 1
$wsInfer.Range("G29").Value = 1
$wsInfer.Range("I29").Value = 0.2592724864350675
$wsInfer.Range("J29").Value = 0.2328233516916538
$wsInfer.Range("K29").Value = 0.3080616184861621
$wsInfer.Range("L29").Value = 0.2457251662216493
